$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure date-like text fields (Y/AA = Startdatum/Slutdatum) remain plain text
# and are not reinterpreted as Excel date serials when (re)written.
foreach ($addr in @("Y106", "AA106", "Y107", "AA107", "Y109", "AA109", "Y110", "AA110", "Y111", "AA111", "Y112", "AA112", "Y113", "AA113", "Y114", "AA114", "Y125", "AA125", "Y126", "AA126", "Y127", "AA127", "Y128", "AA128")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 106 <= content of former row 107
$ws.Range("A106").Value = 130964642
$ws.Range("B106").Value = 99036
$ws.Range("D106").Value = "LC"
$ws.Range("E106").Value = 221952
$ws.Range("F106").Value = "Spindelblomster"
$ws.Range("G106").Value = "Neottia cordata"
$ws.Range("H106").Value = "(L.) Rich."
$ws.Range("Q106").Value = 509917
$ws.Range("R106").Value = 6719042
$ws.Range("Y106").Value = "2025-07-02"
$ws.Range("AA106").Value = "2025-07-02"
$ws.Range("AC106").Value = "Måttliga förekomster . inventering åt vasa vind"

# Row 107 <= content of former row 106
$ws.Range("A107").Value = 130964537
$ws.Range("B107").Value = 79243
$ws.Range("D107").Value = "NT"
$ws.Range("E107").Value = 6425
$ws.Range("F107").Value = "Garnlav"
$ws.Range("G107").Value = "Alectoria sarmentosa"
$ws.Range("H107").Value = "(Ach.) Ach."
$ws.Range("Q107").Value = 509822
$ws.Range("R107").Value = 6718960
$ws.Range("Y107").Value = "2025-07-02"
$ws.Range("AA107").Value = "2025-07-02"
$ws.Range("AC107").Value = "Rikligt . inventering åt vasa vind"

# Row 109 <= content of former row 111
$ws.Range("A109").Value = 130964641
$ws.Range("B109").Value = 98930
$ws.Range("D109").Value = "LC"
$ws.Range("E109").Value = 219790
$ws.Range("F109").Value = "Fläcknycklar"
$ws.Range("G109").Value = "Dactylorhiza maculata"
$ws.Range("H109").Value = "(L.) Soó"
$ws.Range("Q109").Value = 509932
$ws.Range("R109").Value = 6719045
$ws.Range("Y109").Value = "2025-07-02"
$ws.Range("AA109").Value = "2025-07-02"
$ws.Range("AC109").Value = "Måttlig förekomst . inventering åt vasa vind"

# Row 110 <= content of former row 109
$ws.Range("A110").Value = 130964647
$ws.Range("B110").Value = 92106
$ws.Range("D110").Value = "NT"
$ws.Range("E110").Value = 658
$ws.Range("F110").Value = "Rosenticka"
$ws.Range("G110").Value = "Fomitopsis rosea"
$ws.Range("H110").Value = "(Alb. & Schwein.:Fr.) P.Karst."
$ws.Range("Q110").Value = 509741
$ws.Range("R110").Value = 6718998
$ws.Range("Y110").Value = "2025-07-02"
$ws.Range("AA110").Value = "2025-07-02"
$ws.Range("AC110").Value = "Måttliga förekomster . inventering åt vasa vind"

# Row 111 <= content of former row 110
$ws.Range("A111").Value = 130964544
$ws.Range("B111").Value = 57073
$ws.Range("D111").Value = "LC"
$ws.Range("E111").Value = 100138
$ws.Range("F111").Value = "Tjäder"
$ws.Range("G111").Value = "Tetrao urogallus"
$ws.Range("H111").Value = "Linnaeus, 1758"
$ws.Range("Q111").Value = 509543
$ws.Range("R111").Value = 6718926
$ws.Range("Y111").Value = "2025-07-02"
$ws.Range("AA111").Value = "2025-07-02"
$ws.Range("AC111").Value = "Spillning . inventering åt vasa vind"

# Row 112 <= content of former row 114
$ws.Range("A112").Value = 130964650
$ws.Range("B112").Value = 92267
$ws.Range("D112").Value = "VU"
$ws.Range("E112").Value = 1209
$ws.Range("F112").Value = "Rynkskinn"
$ws.Range("G112").Value = "Hermanssonia centrifuga"
$ws.Range("H112").Value = "(P. Karst.) Zmitr."
$ws.Range("Q112").Value = 509694
$ws.Range("R112").Value = 6718936
$ws.Range("Y112").Value = "2025-07-02"
$ws.Range("AA112").Value = "2025-07-02"
$ws.Range("AC112").Value = "Måttliga förekomster . inventering åt vasa vind"

# Row 113 <= content of former row 112
$ws.Range("A113").Value = 130964533
$ws.Range("B113").Value = 79243
$ws.Range("D113").Value = "NT"
$ws.Range("E113").Value = 6425
$ws.Range("F113").Value = "Garnlav"
$ws.Range("G113").Value = "Alectoria sarmentosa"
$ws.Range("H113").Value = "(Ach.) Ach."
$ws.Range("Q113").Value = 509984
$ws.Range("R113").Value = 6719028
$ws.Range("Y113").Value = "2025-07-02"
$ws.Range("AA113").Value = "2025-07-02"
$ws.Range("AC113").Value = "Rikligt . inventering åt vasa vind"

# Row 114 <= content of former row 113
$ws.Range("A114").Value = 130964645
$ws.Range("B114").Value = 99036
$ws.Range("D114").Value = "LC"
$ws.Range("E114").Value = 221952
$ws.Range("F114").Value = "Spindelblomster"
$ws.Range("G114").Value = "Neottia cordata"
$ws.Range("H114").Value = "(L.) Rich."
$ws.Range("Q114").Value = 509804
$ws.Range("R114").Value = 6719024
$ws.Range("Y114").Value = "2025-07-02"
$ws.Range("AA114").Value = "2025-07-02"
$ws.Range("AC114").Value = "Måttliga förekomster . inventering åt vasa vind"

# Row 125 <= content of former row 126
$ws.Range("A125").Value = 130964390
$ws.Range("B125").Value = 99013
$ws.Range("D125").Value = "VU"
$ws.Range("E125").Value = 220787
$ws.Range("F125").Value = "Knärot"
$ws.Range("G125").Value = "Goodyera repens"
$ws.Range("H125").Value = "(L.) R. Br."
$ws.Range("Q125").Value = 509475
$ws.Range("R125").Value = 6718881
$ws.Range("Y125").Value = "2025-07-03"
$ws.Range("AA125").Value = "2025-07-03"
$ws.Range("AC125").Value = "Måttliga förekomster, Ca 10-15 plantor . inventering åt vasa vind"

# Row 126 <= content of former row 125
$ws.Range("A126").Value = 130964643
$ws.Range("B126").Value = 98930
$ws.Range("D126").Value = "LC"
$ws.Range("E126").Value = 219790
$ws.Range("F126").Value = "Fläcknycklar"
$ws.Range("G126").Value = "Dactylorhiza maculata"
$ws.Range("H126").Value = "(L.) Soó"
$ws.Range("Q126").Value = 509829
$ws.Range("R126").Value = 6719000
$ws.Range("Y126").Value = "2025-07-02"
$ws.Range("AA126").Value = "2025-07-02"
$ws.Range("AC126").Value = "Måttlig förekomst . inventering åt vasa vind"

# Row 127 <= content of former row 128
$ws.Range("A127").Value = 130964546
$ws.Range("B127").Value = 92503
$ws.Range("D127").Value = "VU"
$ws.Range("E127").Value = 898
$ws.Range("F127").Value = "Blackticka"
$ws.Range("G127").Value = "Steccherinum collabens"
$ws.Range("H127").Value = "(Fr.) Vesterholt"
$ws.Range("Q127").Value = 509515
$ws.Range("R127").Value = 6718886
$ws.Range("Y127").Value = "2025-07-02"
$ws.Range("AA127").Value = "2025-07-02"
$ws.Range("AC127").Value = "Betydande förekomst . inventering åt vasa vind"

# Row 128 <= content of former row 127
$ws.Range("A128").Value = 130964538
$ws.Range("B128").Value = 79243
$ws.Range("D128").Value = "NT"
$ws.Range("E128").Value = 6425
$ws.Range("F128").Value = "Garnlav"
$ws.Range("G128").Value = "Alectoria sarmentosa"
$ws.Range("H128").Value = "(Ach.) Ach."
$ws.Range("Q128").Value = 509875
$ws.Range("R128").Value = 6719025
$ws.Range("Y128").Value = "2025-07-02"
$ws.Range("AA128").Value = "2025-07-02"
$ws.Range("AC128").Value = "Enstaka . inventering åt vasa vind"
